$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before the existing "K" column ("gestionnaires_additionnels"),
# shifting it (and anything after it) one column to the right, and giving the
# new column the plain/default cell style (not the style shared by the rest
# of the header row).
$ws.Columns("K:K").Insert()

$c = $ws.Cells.Item(1, 11)
$c.Value = "administration_tutelle"

# Reset formatting on the new header cell so it lands on the workbook's
# default (unapplied) cell style rather than inheriting the style used by
# the column it was inserted in front of.
$c.Font.Name = "Arial"
$c.Font.Size = 10
$c.NumberFormat = "General"
$c.HorizontalAlignment = -4142
$c.VerticalAlignment = -4107
$c.WrapText = $false
$c.Locked = $true
$c.FormulaHidden = $false

# Match the width of the new column to the target layout.
$ws.Columns("K:K").ColumnWidth = 26.76

# Update selection like the original edit (scrolled/selected near the new column).
$ws.Range("B1").Select() | Out-Null
$ws.Range("K2").Select() | Out-Null
